$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells keep their text representation exactly (avoid
# Excel auto-converting numeric-looking strings into floating point numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.769.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.034.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.607"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.10"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.334.93"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.771"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.026.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.698.76"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.92"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0819"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.16"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.11"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.71"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.21"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.65%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.48"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.51"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.41"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.38%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.95"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.530.15"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.05"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0906"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.13"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.94"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.224.04"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.95%  "
